$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 10 through 39 (old expanded rows no longer needed)
$ws.Range("A10:A39").EntireRow.Delete() | Out-Null

# Set the new condensed values for rows 2 through 9
$ws.Range("A2").Value = "('Blue Elemental Blast', ['{U}', 'Instant', 'Choose one —', '• Counter target red spell.', '• Destroy target red permanent.'])"
$ws.Range("A3").Value = "('Brainstorm', ['{U}', 'Instant', 'Draw three cards, then put two cards from your hand on top of your library in any order.'])"
$ws.Range("A4").Value = "('Counterspell', ['{U}{U}', 'Instant', 'Counter target spell.'])"
$ws.Range("A5").Value = "('Gifts Ungiven', ['{3}{U}', 'Instant', 'Search your library for up to four cards with different names and reveal them. Target opponent chooses two of those cards. Put the chosen cards into your graveyard and the rest into your hand. Then shuffle your library.'])"
$ws.Range("A6").Value = "('Jace Beleren', ['{1}{U}{U}', 'Legendary Planeswalker — Jace', '+2: Each player draws a card.', '−1: Target player draws a card.', '−10: Target player mills twenty cards.', 'Loyalty: 3'])"
$ws.Range("A7").Value = "('Mystical Tutor', ['{U}', 'Instant', 'Search your library for an instant or sorcery card and reveal that card. Shuffle your library, then put the card on top of it.'])"
$ws.Range("A8").Value = "('Negate', ['{1}{U}', 'Instant', 'Counter target noncreature spell.'])"
$ws.Range("A9").Value = "('Threads of Disloyalty', ['{1}{U}{U}', 'Enchantment — Aura', 'Enchant creature with converted mana cost 2 or less', 'You control enchanted creature.'])"
